$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 21:22"

# --- Reorder Sudafrica / Egipto rows (swap their text labels) ---
# Before: A53 = Egipto, A54 = Sudafrica
# After:  A53 = Sudafrica, A54 = Egipto
$ws.Range("A53").Value = "Sudafrica"
$ws.Range("A54").Value = "Egipto"

# --- Numeric data updates (B=Casos totales, C=Nuevos casos, D=Casos activos,
#     E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes) ---

# Row 15: Brasil
$ws.Range("B15").Value = 36599
$ws.Range("C15").Value = 2917
$ws.Range("E15").Value = 20226
$ws.Range("G15").Value = 206
$ws.Range("H15").Value = 2347

# Row 18: Suiza
$ws.Range("E18").Value = 8936
$ws.Range("G18").Value = 41
$ws.Range("H18").Value = 1368

# Row 35: Noruega
$ws.Range("E35").Value = 6840
$ws.Range("G35").Value = 3
$ws.Range("H35").Value = 164

# Row 53: now Sudafrica (updated counts)
$ws.Range("B53").Value = 3034
$ws.Range("C53").Value = 251
$ws.Range("D53").Value = 903
$ws.Range("E53").Value = 2081
$ws.Range("F53").Value = 7
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 50

# Row 54: now Egipto (old Egipto counts)
$ws.Range("B54").Value = 3032
$ws.Range("C54").Value = 188
$ws.Range("D54").Value = 701
$ws.Range("E54").Value = 2107
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 19
$ws.Range("H54").Value = 224

# Row 64: Barein
$ws.Range("D64").Value = 755
$ws.Range("E64").Value = 1011

# Row 110: Georgia
$ws.Range("B110").Value = 388
$ws.Range("C110").Value = 18
$ws.Range("D110").Value = 86
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 4

# Row 163: Eritrea
$ws.Range("D163").Value = 3
$ws.Range("E163").Value = 36
